$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 730.7
$ws.Range("I19").Value = 994.5
$ws.Range("J19").Value = 466.9
$ws.Range("K19").Value = 994.5
$ws.Range("L19").Value = 466.9
$ws.Range("M19").Value = -819.5
$ws.Range("N19").Value = -816.9
$ws.Range("H80").Value = 920.0244
$ws.Range("I80").Value = 425.72223
$ws.Range("K80").Value = 1277.16669
$ws.Range("M80").Value = -279.16669
$ws.Range("H81").Value = 40328
$ws.Range("J81").Value = 40328
$ws.Range("L81").Value = 40328
$ws.Range("N81").Value = -42324
$ws.Range("H83").Value = 920.0244
$ws.Range("I83").Value = 425.72223
$ws.Range("K83").Value = 3831.50007
$ws.Range("M83").Value = 1160.49993
$ws.Range("H84").Value = 40328
$ws.Range("J84").Value = 40328
$ws.Range("L84").Value = 120984
$ws.Range("N84").Value = -130968
$ws.Range("H88").Value = 2028.8055
$ws.Range("J88").Value = 1064.5358
$ws.Range("L88").Value = 1064.5358
$ws.Range("N88").Value = -1876.5358
$ws.Range("H91").Value = 2028.8055
$ws.Range("J91").Value = 1064.5358
$ws.Range("L91").Value = 1064.5358
$ws.Range("N91").Value = -3872.5358
$ws.Range("H113").Value = 3511.6667
$ws.Range("J113").Value = 4000
$ws.Range("L113").Value = 4000
$ws.Range("N113").Value = -10508
$ws.Range("H132").Value = 5864.6587
$ws.Range("I132").Value = 6473.9653
$ws.Range("J132").Value = 4392.1665
$ws.Range("K132").Value = 19421.8959
$ws.Range("L132").Value = 13176.4995
$ws.Range("M132").Value = -16891.8959
$ws.Range("N132").Value = -18236.4995
$ws.Range("H137").Value = 1100.8857
$ws.Range("I137").Value = 857.2174
$ws.Range("J137").Value = 1567.9166
$ws.Range("K137").Value = 2571.6522
$ws.Range("L137").Value = 4703.7498
$ws.Range("M137").Value = -21.65219999999999
$ws.Range("N137").Value = -9803.7498

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 125101230
$ws.Range("J88").Value = 142972590
$ws.Range("L88").Value = 142972590
$ws.Range("N88").Value = -142973402
$ws.Range("H91").Value = 125101230
$ws.Range("J91").Value = 142972590
$ws.Range("L91").Value = 142972590
$ws.Range("N91").Value = -142975398

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1500
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 1500
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 1500
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -1846
$ws.Range("H86").Value = 3176375
$ws.Range("I86").Value = 4257116
$ws.Range("J86").Value = 1697.9375
$ws.Range("K86").Value = 4257116
$ws.Range("L86").Value = 1697.9375
$ws.Range("M86").Value = -4255993
$ws.Range("N86").Value = -3943.9375
$ws.Range("H89").Value = 3176375
$ws.Range("I89").Value = 4257116
$ws.Range("J89").Value = 1697.9375
$ws.Range("K89").Value = 21285580
$ws.Range("L89").Value = 8489.6875
$ws.Range("M89").Value = -21279964
$ws.Range("N89").Value = -19721.6875
$ws.Range("H105").Value = 2203
$ws.Range("I105").Value = 1738.8889
$ws.Range("K105").Value = 1738.8889
$ws.Range("M105").Value = 8.111100000000079

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 6507.8887
$ws.Range("I99").Value = 10402.75
$ws.Range("J99").Value = 3392
$ws.Range("K99").Value = 10402.75
$ws.Range("L99").Value = 3392
$ws.Range("M99").Value = -8904.75
$ws.Range("N99").Value = -6388
$ws.Range("H122").Value = 982.4
$ws.Range("I122").Value = 1006
$ws.Range("J122").Value = 966.6667
$ws.Range("K122").Value = 3018
$ws.Range("L122").Value = 2900.0001
$ws.Range("M122").Value = -568
$ws.Range("N122").Value = -7800.0001
$ws.Range("H126").Value = 6507.8887
$ws.Range("I126").Value = 10402.75
$ws.Range("J126").Value = 3392
$ws.Range("K126").Value = 31208.25
$ws.Range("L126").Value = 10176
$ws.Range("M126").Value = -28738.25
$ws.Range("N126").Value = -15116
$ws.Range("H132").Value = 1088.2759
$ws.Range("I132").Value = 1025.3846
$ws.Range("J132").Value = 1633.3334
$ws.Range("K132").Value = 3076.1538
$ws.Range("L132").Value = 4900.0002
$ws.Range("M132").Value = -546.1538
$ws.Range("N132").Value = -9960.0002
$ws.Range("H135").Value = 18631222
$ws.Range("J135").Value = 18631222
$ws.Range("L135").Value = 18631222
$ws.Range("N135").Value = -18641362

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1077
$ws.Range("I5").Value = 910
$ws.Range("K5").Value = 2730
$ws.Range("M5").Value = -2618
$ws.Range("H38").Value = 333.8421
$ws.Range("J38").Value = 71.22221999999999
$ws.Range("L38").Value = 213.66666
$ws.Range("N38").Value = -907.66666
$ws.Range("H107").Value = 173.35294
$ws.Range("J107").Value = 184.4
$ws.Range("L107").Value = 553.2
$ws.Range("N107").Value = -4393.2
$ws.Range("H121").Value = 11806
$ws.Range("I121").Value = 33553.332
$ws.Range("J121").Value = 7968.2354
$ws.Range("K121").Value = 100659.996
$ws.Range("L121").Value = 23904.7062
$ws.Range("M121").Value = -99349.99600000001
$ws.Range("N121").Value = -26524.7062
$ws.Range("H131").Value = 840.4595
$ws.Range("I131").Value = 493.89474
$ws.Range("J131").Value = 1206.2778
$ws.Range("K131").Value = 1481.68422
$ws.Range("L131").Value = 3618.8334
$ws.Range("M131").Value = 3558.31578
$ws.Range("N131").Value = -13698.8334
$ws.Range("H132").Value = 1033.5161
$ws.Range("J132").Value = 1560
$ws.Range("L132").Value = 14040
$ws.Range("N132").Value = -19100
$ws.Range("H135").Value = 1077
$ws.Range("I135").Value = 910
$ws.Range("K135").Value = 8190
$ws.Range("M135").Value = -5655
$ws.Range("H140").Value = 1725.3572
$ws.Range("I140").Value = 1606.875
$ws.Range("J140").Value = 1883.3334
$ws.Range("K140").Value = 4820.625
$ws.Range("L140").Value = 5650.0002
$ws.Range("M140").Value = 359.375
$ws.Range("N140").Value = -16010.0002

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 220.6875
$ws.Range("I2").Value = 141.55556
$ws.Range("J2").Value = 322.42856
$ws.Range("K2").Value = 141.55556
$ws.Range("L2").Value = 322.42856
$ws.Range("M2").Value = -28.55556000000001
$ws.Range("N2").Value = -548.4285600000001
$ws.Range("H122").Value = 1944.8125
$ws.Range("I122").Value = 2124.7778
$ws.Range("J122").Value = 1713.4286
$ws.Range("K122").Value = 6374.3334
$ws.Range("L122").Value = 5140.2858
$ws.Range("M122").Value = -3924.3334
$ws.Range("N122").Value = -10040.2858

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2239.6667
$ws.Range("I40").Value = 1931.2858
$ws.Range("J40").Value = 2856.4285
$ws.Range("K40").Value = 1931.2858
$ws.Range("L40").Value = 2856.4285
$ws.Range("M40").Value = -1795.2858
$ws.Range("N40").Value = -3128.4285
$ws.Range("H122").Value = 5632.3887
$ws.Range("I122").Value = 7179.4
$ws.Range("J122").Value = 3698.625
$ws.Range("K122").Value = 21538.2
$ws.Range("L122").Value = 11095.875
$ws.Range("M122").Value = -19088.2
$ws.Range("N122").Value = -15995.875
$ws.Range("H132").Value = 12507844
$ws.Range("I132").Value = 18528490
$ws.Range("J132").Value = 3423
$ws.Range("K132").Value = 55585470
$ws.Range("L132").Value = 10269
$ws.Range("M132").Value = -55582940
$ws.Range("N132").Value = -15329
$ws.Range("H136").Value = 7965.591
$ws.Range("I136").Value = 8826.5
$ws.Range("J136").Value = 5669.8335
$ws.Range("K136").Value = 26479.5
$ws.Range("L136").Value = 17009.5005
$ws.Range("M136").Value = -23929.5
$ws.Range("N136").Value = -22109.5005

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4711.893
$ws.Range("I136").Value = 1262.8636
$ws.Range("J136").Value = 17358.334
$ws.Range("K136").Value = 3788.5908
$ws.Range("L136").Value = 52075.00199999999
$ws.Range("M136").Value = -1238.5908
$ws.Range("N136").Value = -57175.00199999999
